$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.01739057681450047
$ws.Range("C2").Value = 0.4588917497650545
$ws.Range("D2").Value = 0.3571886251912578
$ws.Range("E2").Value = 0.5976525957370701
$ws.Range("F2").Value = 0.6076122229502987

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2020935392784192
$ws.Range("C3").Value = 0.6184300652536038
$ws.Range("D3").Value = 0.581868977188309
$ws.Range("E3").Value = 0.7628033673157906
$ws.Range("F3").Value = 0.7394066550299295
$ws.Range("G3").Value = 96

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1576667027323599
$ws.Range("C4").Value = 0.7100091307271836
$ws.Range("D4").Value = 0.5735266738964355
$ws.Range("E4").Value = 0.757315438834067
$ws.Range("F4").Value = 0.7489062053468125
$ws.Range("G4").Value = 46
